# Finally fixed process for repairing vacant building stocks; created
# importEntities to replace importOwnerHouseholds etc.
#
# - forsale_stock!A4: rename owner "Lucas" -> "Lucius"
# - forrent_stock: insert Name / Savings / Owner Insurance columns at the
#   front (matching the owners/forsale_stock layout), and populate the new
#   owner data (Butch, Harvey, Lee, Carmine). Row 4's occupancy type is
#   corrected from "Multi Family Dwelling" to "Mobile Home" (the old value
#   is retired from the workbook entirely).
# - forrent_stock becomes the active / selected sheet (was forsale_stock).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# forrent_stock: insert 3 columns at the front, then fill them in.
# ---------------------------------------------------------------------
$wsRent = $wb.Worksheets.Item("forrent_stock")
$wsRent.Range("A1:C1").EntireColumn.Insert() | Out-Null

$wsRent.Range("A1").Value = "Name"
$wsRent.Range("B1").Value = "Savings"
$wsRent.Range("C1").Value = "Owner Insurance"

$wsRent.Range("A2").Value = "Butch"
$wsRent.Range("B2").Value = 100000000
$wsRent.Range("C2").Value = 1

$wsRent.Range("A3").Value = "Harvey"
$wsRent.Range("B3").Value = 100000000
$wsRent.Range("C3").Value = 1

# ---------------------------------------------------------------------
# forsale_stock: rename owner "Lucas" to "Lucius".
# ---------------------------------------------------------------------
$wsForSale = $wb.Worksheets.Item("forsale_stock")
$wsForSale.Range("A4").Value = "Lucius"

# ---------------------------------------------------------------------
# back to forrent_stock: finish filling rows 4-5, fix row 4's occupancy.
# ---------------------------------------------------------------------
$wsRent.Range("A4").Value = "Lee"
$wsRent.Range("B4").Value = 100000000
$wsRent.Range("C4").Value = 1
$wsRent.Range("E4").Value = "Mobile Home"

$wsRent.Range("A5").Value = "Carmine"
$wsRent.Range("B5").Value = 100000000
$wsRent.Range("C5").Value = 1

# ---------------------------------------------------------------------
# sheet selection / active-tab bookkeeping.
# ---------------------------------------------------------------------
$wsForSale.Activate() | Out-Null
$wsForSale.Range("H24").Select() | Out-Null

$wsRent.Activate() | Out-Null
$wsRent.Range("H13").Select() | Out-Null
